# Update "想去人数" (interested-people count) figures that were refreshed
# by the scraper run recorded in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 940
$wsExhibit.Range("F6").Value  = 4888
$wsExhibit.Range("F7").Value  = 378
$wsExhibit.Range("F8").Value  = 555
$wsExhibit.Range("F13").Value = 535
$wsExhibit.Range("F16").Value = 1578
$wsExhibit.Range("F18").Value = 674
$wsExhibit.Range("F21").Value = 245
$wsExhibit.Range("F24").Value = 1035
$wsExhibit.Range("F26").Value = 521
$wsExhibit.Range("F27").Value = 1426

# Sheet "演出"
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 128
$wsShow.Range("F6").Value = 95

# Sheet "全部类型" (aggregated view of all the above, offset by one row)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 940
$wsAll.Range("F8").Value  = 4888
$wsAll.Range("F9").Value  = 378
$wsAll.Range("F10").Value = 555
$wsAll.Range("F12").Value = 128
$wsAll.Range("F16").Value = 95
$wsAll.Range("F19").Value = 535
$wsAll.Range("F23").Value = 1578
$wsAll.Range("F25").Value = 674
$wsAll.Range("F28").Value = 245
$wsAll.Range("F32").Value = 1035
$wsAll.Range("F33").Value = 521
$wsAll.Range("F34").Value = 1426
